$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "L1cam"
$ws.Cells.Item(2, 3).Value = "Erbb2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.861952666666666
$ws.Cells.Item(2, 8).Value = 14.585858
$ws.Cells.Item(2, 9).Value = 0.3995648519435639
$ws.Cells.Item(2, 10).Value = 0.3995648519435638
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.020000333333333
$ws.Cells.Item(2, 14).Value = 9.060001
$ws.Cells.Item(2, 15).Value = 0.291481777372034
$ws.Cells.Item(2, 16).Value = 0.291481777372034
$ws.Cells.Item(2, 17).Value = 14.68309867398422
$ws.Cells.Item(2, 18).Value = 132.147888065858
$ws.Cells.Item(2, 19).Value = 0.1164658732199036
$ws.Cells.Item(2, 20).Value = 0.1164658732199036

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "L1cam"
$ws.Cells.Item(3, 3).Value = "Erbb2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.861952666666666
$ws.Cells.Item(3, 8).Value = 14.585858
$ws.Cells.Item(3, 9).Value = 0.3995648519435639
$ws.Cells.Item(3, 10).Value = 0.3995648519435638
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.076388666666666
$ws.Cells.Item(3, 14).Value = 12.229166
$ws.Cells.Item(3, 15).Value = 0.3934413518781783
$ws.Cells.Item(3, 16).Value = 0.3934413518781784
$ws.Cells.Item(3, 17).Value = 19.81920874826978
$ws.Cells.Item(3, 18).Value = 178.372878734428
$ws.Cells.Item(3, 19).Value = 0.15720533551168
$ws.Cells.Item(3, 20).Value = 0.15720533551168

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "L1cam"
$ws.Cells.Item(4, 3).Value = "Erbb2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.861952666666666
$ws.Cells.Item(4, 8).Value = 14.585858
$ws.Cells.Item(4, 9).Value = 0.3995648519435639
$ws.Cells.Item(4, 10).Value = 0.3995648519435638
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.229698
$ws.Cells.Item(4, 14).Value = 9.689094
$ws.Cells.Item(4, 15).Value = 0.311721195201271
$ws.Cells.Item(4, 16).Value = 0.3117211952012711
$ws.Cells.Item(4, 17).Value = 15.702638803628
$ws.Cells.Item(4, 18).Value = 141.323749232652
$ws.Cells.Item(4, 19).Value = 0.1245528332082666
$ws.Cells.Item(4, 20).Value = 0.1245528332082666

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "L1cam"
$ws.Cells.Item(5, 3).Value = "Erbb2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.861952666666666
$ws.Cells.Item(5, 8).Value = 14.585858
$ws.Cells.Item(5, 9).Value = 0.3995648519435639
$ws.Cells.Item(5, 10).Value = 0.3995648519435638
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.03476766666666667
$ws.Cells.Item(5, 14).Value = 0.104303
$ws.Cells.Item(5, 15).Value = 0.003355675548516525
$ws.Cells.Item(5, 16).Value = 0.003355675548516525
$ws.Cells.Item(5, 17).Value = 0.1690387496637778
$ws.Cells.Item(5, 18).Value = 1.521348746974
$ws.Cells.Item(5, 19).Value = 0.001340810003713643
$ws.Cells.Item(5, 20).Value = 0.001340810003713643

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "L1cam"
$ws.Cells.Item(6, 3).Value = "Erbb2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.6001993333333334
$ws.Cells.Item(6, 8).Value = 1.800598
$ws.Cells.Item(6, 9).Value = 0.04932556406896855
$ws.Cells.Item(6, 10).Value = 0.04932556406896854
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.020000333333333
$ws.Cells.Item(6, 14).Value = 9.060001
$ws.Cells.Item(6, 15).Value = 0.291481777372034
$ws.Cells.Item(6, 16).Value = 0.291481777372034
$ws.Cells.Item(6, 17).Value = 1.812602186733111
$ws.Cells.Item(6, 18).Value = 16.313419680598
$ws.Cells.Item(6, 19).Value = 0.01437750308470109
$ws.Cells.Item(6, 20).Value = 0.01437750308470109

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "L1cam"
$ws.Cells.Item(7, 3).Value = "Erbb2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.6001993333333334
$ws.Cells.Item(7, 8).Value = 1.800598
$ws.Cells.Item(7, 9).Value = 0.04932556406896855
$ws.Cells.Item(7, 10).Value = 0.04932556406896854
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.076388666666666
$ws.Cells.Item(7, 14).Value = 12.229166
$ws.Cells.Item(7, 15).Value = 0.3934413518781783
$ws.Cells.Item(7, 16).Value = 0.3934413518781784
$ws.Cells.Item(7, 17).Value = 2.446645760140889
$ws.Cells.Item(7, 18).Value = 22.019811841268
$ws.Cells.Item(7, 19).Value = 0.01940671660944868
$ws.Cells.Item(7, 20).Value = 0.01940671660944868

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "L1cam"
$ws.Cells.Item(8, 3).Value = "Erbb2"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.6001993333333334
$ws.Cells.Item(8, 8).Value = 1.800598
$ws.Cells.Item(8, 9).Value = 0.04932556406896855
$ws.Cells.Item(8, 10).Value = 0.04932556406896854
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.229698
$ws.Cells.Item(8, 14).Value = 9.689094
$ws.Cells.Item(8, 15).Value = 0.311721195201271
$ws.Cells.Item(8, 16).Value = 0.3117211952012711
$ws.Cells.Item(8, 17).Value = 1.938462586468001
$ws.Cells.Item(8, 18).Value = 17.446163278212
$ws.Cells.Item(8, 19).Value = 0.01537582378555574
$ws.Cells.Item(8, 20).Value = 0.01537582378555574

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "L1cam"
$ws.Cells.Item(9, 3).Value = "Erbb2"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.6001993333333334
$ws.Cells.Item(9, 8).Value = 1.800598
$ws.Cells.Item(9, 9).Value = 0.04932556406896855
$ws.Cells.Item(9, 10).Value = 0.04932556406896854
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.03476766666666667
$ws.Cells.Item(9, 14).Value = 0.104303
$ws.Cells.Item(9, 15).Value = 0.003355675548516525
$ws.Cells.Item(9, 16).Value = 0.003355675548516525
$ws.Cells.Item(9, 17).Value = 0.02086753035488889
$ws.Cells.Item(9, 18).Value = 0.187807773194
$ws.Cells.Item(9, 19).Value = 0.000165520589263023
$ws.Cells.Item(9, 20).Value = 0.000165520589263023

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "L1cam"
$ws.Cells.Item(10, 3).Value = "Erbb2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.206754333333333
$ws.Cells.Item(10, 8).Value = 12.620263
$ws.Cells.Item(10, 9).Value = 0.3457193616641432
$ws.Cells.Item(10, 10).Value = 0.3457193616641432
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.020000333333333
$ws.Cells.Item(10, 14).Value = 9.060001
$ws.Cells.Item(10, 15).Value = 0.291481777372034
$ws.Cells.Item(10, 16).Value = 0.291481777372034
$ws.Cells.Item(10, 17).Value = 12.70439948891811
$ws.Cells.Item(10, 18).Value = 114.339595400263
$ws.Cells.Item(10, 19).Value = 0.1007708940097895
$ws.Cells.Item(10, 20).Value = 0.1007708940097895

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "L1cam"
$ws.Cells.Item(11, 3).Value = "Erbb2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.206754333333333
$ws.Cells.Item(11, 8).Value = 12.620263
$ws.Cells.Item(11, 9).Value = 0.3457193616641432
$ws.Cells.Item(11, 10).Value = 0.3457193616641432
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 4.076388666666666
$ws.Cells.Item(11, 14).Value = 12.229166
$ws.Cells.Item(11, 15).Value = 0.3934413518781783
$ws.Cells.Item(11, 16).Value = 0.3934413518781784
$ws.Cells.Item(11, 17).Value = 17.14836568785088
$ws.Cells.Item(11, 18).Value = 154.335291190658
$ws.Cells.Item(11, 19).Value = 0.1360202930236014
$ws.Cells.Item(11, 20).Value = 0.1360202930236014

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "L1cam"
$ws.Cells.Item(12, 3).Value = "Erbb2"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.206754333333333
$ws.Cells.Item(12, 8).Value = 12.620263
$ws.Cells.Item(12, 9).Value = 0.3457193616641432
$ws.Cells.Item(12, 10).Value = 0.3457193616641432
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.229698
$ws.Cells.Item(12, 14).Value = 9.689094
$ws.Cells.Item(12, 15).Value = 0.311721195201271
$ws.Cells.Item(12, 16).Value = 0.3117211952012711
$ws.Cells.Item(12, 17).Value = 13.586546056858
$ws.Cells.Item(12, 18).Value = 122.278914511722
$ws.Cells.Item(12, 19).Value = 0.1077680526221672
$ws.Cells.Item(12, 20).Value = 0.1077680526221672

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "L1cam"
$ws.Cells.Item(13, 3).Value = "Erbb2"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.206754333333333
$ws.Cells.Item(13, 8).Value = 12.620263
$ws.Cells.Item(13, 9).Value = 0.3457193616641432
$ws.Cells.Item(13, 10).Value = 0.3457193616641432
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.03476766666666667
$ws.Cells.Item(13, 14).Value = 0.104303
$ws.Cells.Item(13, 15).Value = 0.003355675548516525
$ws.Cells.Item(13, 16).Value = 0.003355675548516525
$ws.Cells.Item(13, 17).Value = 0.1462590324098889
$ws.Cells.Item(13, 18).Value = 1.316331291689
$ws.Cells.Item(13, 19).Value = 0.001160122008585107
$ws.Cells.Item(13, 20).Value = 0.001160122008585107

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "L1cam"
$ws.Cells.Item(14, 3).Value = "Erbb2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2.499212666666667
$ws.Cells.Item(14, 8).Value = 7.497638
$ws.Cells.Item(14, 9).Value = 0.2053902223233243
$ws.Cells.Item(14, 10).Value = 0.2053902223233243
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.020000333333333
$ws.Cells.Item(14, 14).Value = 9.060001
$ws.Cells.Item(14, 15).Value = 0.291481777372034
$ws.Cells.Item(14, 16).Value = 0.291481777372034
$ws.Cells.Item(14, 17).Value = 7.547623086404222
$ws.Cells.Item(14, 18).Value = 67.928607777638
$ws.Cells.Item(14, 19).Value = 0.05986750705763979
$ws.Cells.Item(14, 20).Value = 0.05986750705763979

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "L1cam"
$ws.Cells.Item(15, 3).Value = "Erbb2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2.499212666666667
$ws.Cells.Item(15, 8).Value = 7.497638
$ws.Cells.Item(15, 9).Value = 0.2053902223233243
$ws.Cells.Item(15, 10).Value = 0.2053902223233243
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 4.076388666666666
$ws.Cells.Item(15, 14).Value = 12.229166
$ws.Cells.Item(15, 15).Value = 0.3934413518781783
$ws.Cells.Item(15, 16).Value = 0.3934413518781784
$ws.Cells.Item(15, 17).Value = 10.18776218998978
$ws.Cells.Item(15, 18).Value = 91.68985970990799
$ws.Cells.Item(15, 19).Value = 0.08080900673344832
$ws.Cells.Item(15, 20).Value = 0.08080900673344833

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "L1cam"
$ws.Cells.Item(16, 3).Value = "Erbb2"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.499212666666667
$ws.Cells.Item(16, 8).Value = 7.497638
$ws.Cells.Item(16, 9).Value = 0.2053902223233243
$ws.Cells.Item(16, 10).Value = 0.2053902223233243
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 3.229698
$ws.Cells.Item(16, 14).Value = 9.689094
$ws.Cells.Item(16, 15).Value = 0.311721195201271
$ws.Cells.Item(16, 16).Value = 0.3117211952012711
$ws.Cells.Item(16, 17).Value = 8.071702151108001
$ws.Cells.Item(16, 18).Value = 72.645319359972
$ws.Cells.Item(16, 19).Value = 0.06402448558528143
$ws.Cells.Item(16, 20).Value = 0.06402448558528144

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "L1cam"
$ws.Cells.Item(17, 3).Value = "Erbb2"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 2.499212666666667
$ws.Cells.Item(17, 8).Value = 7.497638
$ws.Cells.Item(17, 9).Value = 0.2053902223233243
$ws.Cells.Item(17, 10).Value = 0.2053902223233243
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.03476766666666667
$ws.Cells.Item(17, 14).Value = 0.104303
$ws.Cells.Item(17, 15).Value = 0.003355675548516525
$ws.Cells.Item(17, 16).Value = 0.003355675548516525
$ws.Cells.Item(17, 17).Value = 0.08689179292377779
$ws.Cells.Item(17, 18).Value = 0.7820261363140001
$ws.Cells.Item(17, 19).Value = 0.0006892229469547523
$ws.Cells.Item(17, 20).Value = 0.0006892229469547524

